$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 231 (existing data from row 231 onward shifts down to row 233+)
$ws.Rows.Item(231).Resize(2).Insert()

# Populate new row 231
$ws.Cells.Item(231,1).Value2  = 1
$ws.Cells.Item(231,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(231,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(231,4).Value2  = 45215
$ws.Cells.Item(231,5).Value2  = 15
$ws.Cells.Item(231,6).Value2  = 100114001
$ws.Cells.Item(231,7).Value2  = "Papa"
$ws.Cells.Item(231,8).Value2  = "Asterix"
$ws.Cells.Item(231,9).Value2  = "1a (cosecha)"
$ws.Cells.Item(231,10).Value2 = 700
$ws.Cells.Item(231,11).Value2 = 31000
$ws.Cells.Item(231,12).Value2 = 32000
$ws.Cells.Item(231,13).Value2 = 31429
$ws.Cells.Item(231,14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(231,15).Value2 = "Región del Maule"
$ws.Cells.Item(231,16).Value2 = 1257
$ws.Cells.Item(231,17).Value2 = 25
$ws.Cells.Item(231,18).Value2 = "Hortaliza"

# Populate new row 232
$ws.Cells.Item(232,1).Value2  = 1
$ws.Cells.Item(232,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(232,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(232,4).Value2  = 45215
$ws.Cells.Item(232,5).Value2  = 15
$ws.Cells.Item(232,6).Value2  = 100114001
$ws.Cells.Item(232,7).Value2  = "Papa"
$ws.Cells.Item(232,8).Value2  = "Red Lady"
$ws.Cells.Item(232,9).Value2  = "1a (cosecha)"
$ws.Cells.Item(232,10).Value2 = 600
$ws.Cells.Item(232,11).Value2 = 29000
$ws.Cells.Item(232,12).Value2 = 30000
$ws.Cells.Item(232,13).Value2 = 29583
$ws.Cells.Item(232,14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(232,15).Value2 = "Región del Maule"
$ws.Cells.Item(232,16).Value2 = 1183
$ws.Cells.Item(232,17).Value2 = 25
$ws.Cells.Item(232,18).Value2 = "Hortaliza"

Write-Host "Done. Dimension:" $ws.UsedRange.Address()
